$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 412, shifting existing rows 412.. down by one
$ws.Rows.Item(412).Insert()

# Populate the newly inserted row 412 with the new data record
$ws.Cells.Item(412, 1).Value = 3
$ws.Cells.Item(412, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(412, 3).Value = "Coquimbo"
$ws.Cells.Item(412, 4).Value = 44946
$ws.Cells.Item(412, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(412, 5).Value = 5
$ws.Cells.Item(412, 6).Value = 100112012
$ws.Cells.Item(412, 7).Value = "Espinaca"
$ws.Cells.Item(412, 8).Value = "Sin especificar"
$ws.Cells.Item(412, 9).Value = "Primera"
$ws.Cells.Item(412, 10).Value = 230
$ws.Cells.Item(412, 11).Value = 4000
$ws.Cells.Item(412, 12).Value = 4200
$ws.Cells.Item(412, 13).Value = 4104
$ws.Cells.Item(412, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(412, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(412, 16).Value = 1368
$ws.Cells.Item(412, 17).Value = 3
$ws.Cells.Item(412, 18).Value = "Hortaliza"
